# Apply the "home-sales" edit:
#   - Column B values were originally recorded in thousands (e.g. 43.2)
#     and are rescaled to full units (e.g. 43200) by multiplying by 1000.
#   - Column D, which only ever held empty, pre-styled placeholder cells,
#     is removed entirely (this also shrinks the sheet dimension from
#     A1:D62 down to A1:C62 and each row's "spans" from 1:4 to 1:3).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count()

# Row 1 is the header ("Value") - leave it alone, only rescale the
# numeric data rows below it.
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $val = $cell.Value()
    if ($val -ne $null) {
        $cell.Value = $val * 1000
    }
}

# Remove column D (shifts nothing else; there is no column E in use).
$ws.Columns.Item(4).Delete()
